# Applies the commit: "added user name info at settings modal and changed a user's info"
# Concretely: on the "User" sheet, row 15 (user_id 1236) has its username,
# first_name and last_name updated from DENSU/DENY/SUMARGA to DION/DION/WIYOYO,
# and the sheet's active-cell selection moves from G16 to G15.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("User")
$ws.Activate()

# Update the user's info on row 15 (user_id = 1236)
$ws.Range("B15").Value = "DION"     # username
$ws.Range("D15").Value = "DION"     # first_name
$ws.Range("E15").Value = "WIYOYO"   # last_name

# Update the saved selection for the User sheet
$ws.Range("G15").Select()
